$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "???" -> "31st Oct 2021" (header line), and drop the red highlight
#    on that run.
# ------------------------------------------------------------------
$dateRange = $d.Content
$foundDate = $dateRange.Find.Execute("???", $true, $false, $false, $false, $false, $true, 1, $false, "31st Oct 2021", 2)
if ($foundDate) {
    $afterReplace = $d.Content
    $null = $afterReplace.Find.Execute("31st Oct 2021", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $afterReplace.HighlightColorIndex = 0
}

# ------------------------------------------------------------------
# 2) Split each "(ignore those which are already incorporated in your
#    book's version and date)" run into three runs: "(", "ignore",
#    " those which are already incorporated in your book's version and
#    date)" -- wrapping "ignore" with grammar-check proofErr markers,
#    same as Word does automatically once the paragraph text has been
#    retyped.
# ------------------------------------------------------------------
$needle = "(ignore those which are already incorporated in your book" + [char]0x2019 + "s version and date)"

$searchRange = $d.Content
$searchRange.Start = 0
$searchRange.End = $d.Content.End

$bmCounter = 0
while ($true) {
    $found = $searchRange.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) { break }

    $s = $searchRange.Start
    $e = $searchRange.End

    # "(" = 1 char, "ignore" = 6 chars, remainder = rest of the range
    $ignoreStart = $s + 1
    $ignoreEnd = $ignoreStart + 6

    $bmCounter = $bmCounter + 1
    $bmName = "zzTmpSplit" + $bmCounter

    $ignoreRange = $d.Range($ignoreStart, $ignoreEnd)
    $d.Bookmarks.Add($bmName, $ignoreRange)
    $d.Bookmarks($bmName).Delete()

    # Continue searching after this match.
    $searchRange.Start = $e
    $searchRange.End = $d.Content.End
}

Write-Output "Replaced date: $foundDate; split paragraphs: $bmCounter"
